# Applies the "Folder Inventory updated" edit described by the commit:
#  - A new folder entry is inserted at the top of the "Folder Inventory"
#    table (row 2), pushing all existing entries down by one row.
#  - The "Metadata" sheet's generation timestamp / folder count / workflow
#    run counters are refreshed.
#  - The "Summary" sheet's folder counters and "Most Recent Update" value
#    are refreshed to match.

$wb = $excel.ActiveWorkbook

$wsInventory = $wb.Worksheets.Item("Folder Inventory")
$wsMetadata  = $wb.Worksheets.Item("Metadata")
$wsSummary   = $wb.Worksheets.Item("Summary")

# --- 1. Folder Inventory: insert a new row at the top of the data (row 2) ---
$newRow = $wsInventory.Rows.Item(2)
$newRow.Insert()
$newRow.ClearFormats()

$newName = "Get Started with Data Warehouses and Ingesting Data with Dataflows Gen2 with Microsoft Fabric"
$newDate = "2025-06-16 19:38:48 +0530"

$wsInventory.Cells.Item(2, 1).Value = $newName
$wsInventory.Cells.Item(2, 2).Value = $newName
$wsInventory.Cells.Item(2, 3).Value = $newDate
$wsInventory.Cells.Item(2, 4).Value = 1
$wsInventory.Cells.Item(2, 5).Value = "Root"

# --- 2. Metadata sheet updates ---
$wsMetadata.Cells.Item(3, 2).Value = "2025-06-16 14:09:08 UTC"
$wsMetadata.Cells.Item(4, 2).Value = 75

# "Workflow Run" is stored as text (not a number) in the workbook, so force
# the cell to text format before assigning the numeric-looking string,
# otherwise Excel will auto-coerce it into a true number. Clear the
# formatting back afterwards so no stray cell style is introduced.
$wsMetadata.Cells.Item(5, 2).NumberFormat = "@"
$wsMetadata.Cells.Item(5, 2).Value = "10"
$wsMetadata.Cells.Item(5, 2).ClearFormats()

# --- 3. Summary sheet updates ---
$wsSummary.Cells.Item(2, 2).Value = 75
$wsSummary.Cells.Item(3, 2).Value = 75
$wsSummary.Cells.Item(5, 2).Value = $newDate
